# Scheduled runner update: refresh cached market-board pricing/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the Asura data-center
# profit sheets. Values below are the freshly pulled/recomputed figures;
# cells that are cleared correspond to rows whose computed profit columns
# collapsed to blank (no data) in the refreshed pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = 0

$ws.Range("H137").Value = 2795.3076
$ws.Range("I137").Value = 1731
$ws.Range("J137").Value = 5190
$ws.Range("K137").Value = 5193
$ws.Range("L137").Value = 15570
$ws.Range("M137").Value = -2643
$ws.Range("N137").Value = -20670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19827.184
$ws.Range("I32").Value = 23225.322
$ws.Range("J32").Value = 4778.2856
$ws.Range("K32").Value = 23225.322
$ws.Range("L32").Value = 4778.2856
$ws.Range("M32").Value = -22938.322
$ws.Range("N32").Value = -5352.2856

$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20626

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0

$ws.Range("H55").Value = 22500

$ws.Range("H74").Value = 1264.5454
$ws.Range("I74").Value = 1212.2222
$ws.Range("K74").Value = 1212.2222
$ws.Range("M74").Value = -338.2221999999999

$ws.Range("H77").Value = 1264.5454
$ws.Range("I77").Value = 1212.2222
$ws.Range("K77").Value = 6061.111
$ws.Range("M77").Value = -1693.111

$ws.Range("H109").Value = 36000
$ws.Range("J109").Value = 36000
$ws.Range("L109").Value = 36000
$ws.Range("N109").Value = -38774

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080

$ws.Range("H132").Value = 3948
$ws.Range("I132").Value = 3375
$ws.Range("J132").Value = 4364.727
$ws.Range("K132").Value = 10125
$ws.Range("L132").Value = 13094.181
$ws.Range("M132").Value = -7595
$ws.Range("N132").Value = -18154.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 29000
$ws.Range("J52").Value = 29000
$ws.Range("L52").Value = 29000
$ws.Range("N52").Value = -29526

$ws.Range("H99").Value = 2268.0908
$ws.Range("I99").Value = 2170
$ws.Range("J99").Value = 2409.7778
$ws.Range("K99").Value = 2170
$ws.Range("L99").Value = 2409.7778
$ws.Range("M99").Value = -672
$ws.Range("N99").Value = -5405.7778

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0

$ws.Range("H121").Value = 29000
$ws.Range("J121").Value = 29000
$ws.Range("L121").Value = 29000
$ws.Range("N121").Value = -32494

$ws.Range("H130").Value = 100000
$ws.Range("J130").Value = 100000
$ws.Range("L130").Value = 100000
$ws.Range("N130").Value = -110040

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 27750
$ws.Range("J109").Value = 27750
$ws.Range("L109").Value = 27750
$ws.Range("N109").Value = -29830

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3885.8572
$ws.Range("I5").Value = 4640.2
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 13920.6
$ws.Range("L5").Value = 6000
$ws.Range("M5").Value = -13808.6
$ws.Range("N5").Value = -6224

$ws.Range("H34").Value = 466.7619
$ws.Range("I34").Value = 224
$ws.Range("J34").Value = 616.1539
$ws.Range("K34").Value = 672
$ws.Range("L34").Value = 1848.4617
$ws.Range("M34").Value = -588
$ws.Range("N34").Value = -2016.4617

$ws.Range("H36").Value = 6999.6665
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2831

$ws.Range("H39").Value = 3027.9092
$ws.Range("J39").Value = 3027.9092
$ws.Range("L39").Value = 9083.7276
$ws.Range("N39").Value = -9671.7276

$ws.Range("H55").Value = 4683.846
$ws.Range("I55").Value = 2122
$ws.Range("J55").Value = 5822.4443
$ws.Range("K55").Value = 6366
$ws.Range("L55").Value = 17467.3329
$ws.Range("M55").Value = -6189
$ws.Range("N55").Value = -17821.3329

$ws.Range("H58").Value = 2400
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").Value = 7200
$ws.Range("N58").Value = -7456

$ws.Range("H129").Value = 2942028.2
$ws.Range("J129").Value = 3334218.8
$ws.Range("L129").Value = 10002656.4
$ws.Range("N129").Value = -10012656.4

$ws.Range("H134").Value = 4385.136
$ws.Range("I134").Value = 2963.3333
$ws.Range("J134").Value = 6091.3
$ws.Range("K134").Value = 8889.999899999999
$ws.Range("L134").Value = 18273.9
$ws.Range("M134").Value = -3819.999899999999
$ws.Range("N134").Value = -28413.9

$ws.Range("H135").Value = 3885.8572
$ws.Range("I135").Value = 4640.2
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 41761.8
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -39226.8
$ws.Range("N135").Value = -23070

$ws.Range("H140").Value = 1984.6061
$ws.Range("I140").Value = 835.4091
$ws.Range("J140").Value = 4283
$ws.Range("K140").Value = 2506.2273
$ws.Range("L140").Value = 12849
$ws.Range("M140").Value = 2673.7727
$ws.Range("N140").Value = -23209

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24579.285
$ws.Range("I57").Value = 4018.3333
$ws.Range("K57").Value = 4018.3333
$ws.Range("M57").Value = -3198.3333

$ws.Range("H123").Value = 21461.684
$ws.Range("J123").Value = 21461.684
$ws.Range("L123").Value = 21461.684
$ws.Range("N123").Value = -26361.684

$ws.Range("H132").Value = 3858.8572
$ws.Range("I132").Value = 3670.6667
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11012.0001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -8482.000100000001
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 400.66666
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 1002
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 1002
$ws.Range("M55").Value = 73
$ws.Range("N55").Value = -1348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11867.5
$ws.Range("I54").Value = 5656.6665
$ws.Range("J54").Value = 30500
$ws.Range("K54").Value = 5656.6665
$ws.Range("L54").Value = 30500
$ws.Range("M54").Value = -5136.6665
$ws.Range("N54").Value = -31540

$ws.Range("H109").Value = 28999
$ws.Range("J109").Value = 28999
$ws.Range("L109").Value = 28999
$ws.Range("N109").Value = -31773

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0

$ws.Range("H113").Value = 698.0909
$ws.Range("I113").Value = 535.3
$ws.Range("K113").Value = 1605.9
$ws.Range("M113").Value = 564.1000000000001

$ws.Range("H118").Value = 28463.455
$ws.Range("J118").Value = 28463.455
$ws.Range("L118").Value = 28463.455
$ws.Range("N118").Value = -31777.455

$ws.Range("H121").Value = 26843.375
$ws.Range("J121").Value = 26843.375
$ws.Range("L121").Value = 26843.375
$ws.Range("N121").Value = -30337.375

$ws.Range("H123").Value = 36112.42
$ws.Range("J123").Value = 36112.42
$ws.Range("L123").Value = 36112.42
$ws.Range("N123").Value = -45912.42

$ws.Range("H139").Value = 62750
$ws.Range("J139").Value = 64571.43
$ws.Range("L139").Value = 64571.43
$ws.Range("N139").Value = -74851.42999999999

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("N141").Value = 0
